$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 585.5
$ws.Range("J12").Value = 448
$ws.Range("L12").Value = 448
$ws.Range("N12").Value = -788
$ws.Range("H64").Value = 5938.3076
$ws.Range("I64").Value = 5950.4
$ws.Range("K64").Value = 5950.4
$ws.Range("M64").Value = -5702.4
$ws.Range("H67").Value = 5938.3076
$ws.Range("I67").Value = 5950.4
$ws.Range("K67").Value = 5950.4
$ws.Range("M67").Value = -5092.4
$ws.Range("H98").Value = 1421.5
$ws.Range("I98").Value = 1421.5
$ws.Range("K98").Value = 1421.5
$ws.Range("M98").Value = 76.5
$ws.Range("H113").Value = 2779.8
$ws.Range("J113").Value = 2912.25
$ws.Range("L113").Value = 2912.25
$ws.Range("N113").Value = -9420.25
$ws.Range("H122").Value = 1421.5
$ws.Range("I122").Value = 1421.5
$ws.Range("K122").Value = 4264.5
$ws.Range("M122").Value = -1814.5
$ws.Range("H125").Value = 2571.4285
$ws.Range("J125").Value = 2624.75
$ws.Range("L125").Value = 23622.75
$ws.Range("N125").Value = -28542.75
$ws.Range("H135").Value = 3849144.5
$ws.Range("J135").Value = 5333.4
$ws.Range("L135").Value = 48000.6
$ws.Range("N135").Value = -53070.6
$ws.Range("H137").Value = 12507756
$ws.Range("I137").Value = 25000932
$ws.Range("J137").Value = 14580.5
$ws.Range("K137").Value = 75002796
$ws.Range("L137").Value = 43741.5
$ws.Range("M137").Value = -75000246
$ws.Range("N137").Value = -48841.5
$ws.Range("H141").Value = 2118.1667
$ws.Range("I141").Value = 1541.8
$ws.Range("K141").Value = 4625.4
$ws.Range("M141").Value = 554.6000000000004

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 125256.555
$ws.Range("I32").Value = 186718
$ws.Range("J32").Value = 2333.6667
$ws.Range("K32").Value = 186718
$ws.Range("L32").Value = 2333.6667
$ws.Range("M32").Value = -186431
$ws.Range("N32").Value = -2907.6667
$ws.Range("H63").Value = 9186.143
$ws.Range("H66").Value = 9186.143
$ws.Range("H74").Value = 2077.652
$ws.Range("I74").Value = 501.05884
$ws.Range("K74").Value = 501.05884
$ws.Range("M74").Value = 372.94116
$ws.Range("H77").Value = 2077.652
$ws.Range("I77").Value = 501.05884
$ws.Range("K77").Value = 2505.2942
$ws.Range("M77").Value = 1862.7058
$ws.Range("H110").Value = 29690000
$ws.Range("I110").Value = 43183590
$ws.Range("K110").Value = 43183590
$ws.Range("M110").Value = -43181545
$ws.Range("H122").Value = 1739.6666
$ws.Range("I122").Value = 1460
$ws.Range("K122").Value = 4380
$ws.Range("M122").Value = -1930
$ws.Range("H132").Value = 4196.1113
$ws.Range("I132").Value = 3117.9546
$ws.Range("K132").Value = 9353.863799999999
$ws.Range("M132").Value = -6823.863799999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 17328.75
$ws.Range("J46").Value = 17328.75
$ws.Range("L46").Value = 17328.75
$ws.Range("N46").Value = -17924.75
$ws.Range("H82").Value = 30191.166
$ws.Range("I82").Value = 16572.273
$ws.Range("K82").Value = 16572.273
$ws.Range("M82").Value = -16189.273
$ws.Range("H85").Value = 30191.166
$ws.Range("I85").Value = 16572.273
$ws.Range("K85").Value = 16572.273
$ws.Range("M85").Value = -15246.273
$ws.Range("H107").Value = 2199.923
$ws.Range("I107").Value = 2093.65
$ws.Range("K107").Value = 2093.65
$ws.Range("M107").Value = -173.6500000000001
$ws.Range("H134").Value = 3164.8928
$ws.Range("I134").Value = 1477.4546
$ws.Range("K134").Value = 4432.3638
$ws.Range("M134").Value = -1897.3638
$ws.Range("H135").Value = 103922
$ws.Range("J135").Value = 103922
$ws.Range("L135").Value = 103922
$ws.Range("N135").Value = -114062

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37041040
$ws.Range("J31").Value = 5044.8335
$ws.Range("L31").Value = 5044.8335
$ws.Range("N31").Value = -5634.8335
$ws.Range("H34").Value = 37041040
$ws.Range("J34").Value = 5044.8335
$ws.Range("L34").Value = 5044.8335
$ws.Range("N34").Value = -5448.8335
$ws.Range("H58").Value = 3715.44
$ws.Range("I58").Value = 2469.8
$ws.Range("K58").Value = 2469.8
$ws.Range("M58").Value = -2266.8
$ws.Range("H99").Value = 3673.8
$ws.Range("J99").Value = 3829
$ws.Range("L99").Value = 3829
$ws.Range("N99").Value = -6825
$ws.Range("H107").Value = 1780.9166
$ws.Range("J107").Value = 2057.25
$ws.Range("L107").Value = 2057.25
$ws.Range("N107").Value = -5897.25
$ws.Range("H126").Value = 3673.8
$ws.Range("J126").Value = 3829
$ws.Range("L126").Value = 11487
$ws.Range("N126").Value = -16427
$ws.Range("H132").Value = 208526
$ws.Range("I132").Value = 7196
$ws.Range("J132").Value = 309191
$ws.Range("K132").Value = 21588
$ws.Range("L132").Value = 927573
$ws.Range("M132").Value = -19058
$ws.Range("N132").Value = -932633
$ws.Range("H136").Value = 3715.44
$ws.Range("I136").Value = 2469.8
$ws.Range("K136").Value = 7409.400000000001
$ws.Range("M136").Value = -4859.400000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 367
$ws.Range("I5").Value = 367
$ws.Range("K5").Value = 1101
$ws.Range("M5").Value = -989
$ws.Range("H38").Value = 74.111115
$ws.Range("I38").Value = 74.111115
$ws.Range("K38").Value = 222.333345
$ws.Range("M38").Value = 124.666655
$ws.Range("H129").Value = 35715508
$ws.Range("J129").Value = 62501332
$ws.Range("L129").Value = 187503996
$ws.Range("N129").Value = -187513996
$ws.Range("H135").Value = 367
$ws.Range("I135").Value = 367
$ws.Range("K135").Value = 3303
$ws.Range("M135").Value = -768

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 640.2143
$ws.Range("I97").Value = 635.61536
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 635.61536
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -139.61536
$ws.Range("N97").Value = -1692
$ws.Range("H132").Value = 4811.931
$ws.Range("I132").Value = 2610.3076
$ws.Range("K132").Value = 7830.9228
$ws.Range("M132").Value = -5300.9228

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 51450.332
$ws.Range("J74").Value = 58158.5
$ws.Range("L74").Value = 58158.5
$ws.Range("N74").Value = -60154.5
$ws.Range("H77").Value = 51450.332
$ws.Range("J77").Value = 58158.5
$ws.Range("L77").Value = 174475.5
$ws.Range("N77").Value = -184459.5
$ws.Range("H122").Value = 7349.6665
$ws.Range("I122").Value = 5310
$ws.Range("K122").Value = 15930
$ws.Range("M122").Value = -13480
$ws.Range("H132").Value = 5967.905
$ws.Range("I132").Value = 5017.0713
$ws.Range("K132").Value = 15051.2139
$ws.Range("M132").Value = -12521.2139
$ws.Range("H136").Value = 3966.6445
$ws.Range("I136").Value = 3091.238
$ws.Range("K136").Value = 9273.714
$ws.Range("M136").Value = -6723.714

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 860.7368
$ws.Range("I113").Value = 661.5
$ws.Range("J113").Value = 1202.2858
$ws.Range("K113").Value = 1984.5
$ws.Range("L113").Value = 3606.8574
$ws.Range("M113").Value = 185.5
$ws.Range("N113").Value = -7946.857400000001
$ws.Range("H126").Value = 3041.2122
$ws.Range("I126").Value = 2045.375
$ws.Range("J126").Value = 5696.778
$ws.Range("K126").Value = 6136.125
$ws.Range("L126").Value = 17090.334
$ws.Range("M126").Value = -3666.125
$ws.Range("N126").Value = -22030.334
$ws.Range("H128").Value = 49465
$ws.Range("J128").Value = 49465
$ws.Range("L128").Value = 49465
$ws.Range("N128").Value = -59425
$ws.Range("H132").Value = 6293.0386
$ws.Range("I132").Value = 3717.4443
$ws.Range("J132").Value = 7656.5884
$ws.Range("K132").Value = 11152.3329
$ws.Range("L132").Value = 22969.7652
$ws.Range("M132").Value = -8622.332900000001
$ws.Range("N132").Value = -28029.7652
$ws.Range("H136").Value = 5515.154
$ws.Range("I136").Value = 2274.75
$ws.Range("K136").Value = 6824.25
$ws.Range("M136").Value = -4274.25
